# Update weekly price records: dates and associated volume/price/origin
# values are reshuffled to their corresponding week across the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) - numeric date serials
$ws.Range("D2").Value  = 44188
$ws.Range("D3").Value  = 44188
$ws.Range("D4").Value  = 44491
$ws.Range("D5").Value  = 44491
$ws.Range("D6").Value  = 44525
$ws.Range("D7").Value  = 44525
$ws.Range("D8").Value  = 44230
$ws.Range("D9").Value  = 44230
$ws.Range("D10").Value = 44308
$ws.Range("D11").Value = 44308
$ws.Range("D12").Value = 44293
$ws.Range("D13").Value = 44293
$ws.Range("D14").Value = 44358
$ws.Range("D15").Value = 44358
$ws.Range("D16").Value = 44328
$ws.Range("D17").Value = 44328
$ws.Range("D18").Value = 44335
$ws.Range("D19").Value = 44335
$ws.Range("D20").Value = 44321
$ws.Range("D21").Value = 44321
$ws.Range("D22").Value = 44554
$ws.Range("D23").Value = 44554

# Column J (Volumen)
$ws.Range("J6").Value  = 200
$ws.Range("J7").Value  = 100
$ws.Range("J8").Value  = 100
$ws.Range("J9").Value  = 50
$ws.Range("J10").Value = 200
$ws.Range("J11").Value = 100
$ws.Range("J14").Value = 200
$ws.Range("J15").Value = 100
$ws.Range("J16").Value = 100
$ws.Range("J17").Value = 50
$ws.Range("J18").Value = 150
$ws.Range("J19").Value = 50

# Column M (Precio promedio ponderado)
$ws.Range("M10").Value = 650
$ws.Range("M18").Value = 633

# Column O (Origen)
$ws.Range("O4").Value  = "Región Metropolitana"
$ws.Range("O5").Value  = "Región Metropolitana"
$ws.Range("O18").Value = "Región de Ñuble"
$ws.Range("O19").Value = "Región de Ñuble"

# Column P (Precio $/Kg)
$ws.Range("P10").Value = 108
$ws.Range("P18").Value = 106
